$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.65
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 2.52
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 2.05
$ws.Range("L5").Value = 3.1
$ws.Range("O5").Value = 1.27
$ws.Range("P5").Value = 3.1
$ws.Range("S5").Value = 1.82
$ws.Range("T5").Value = 1.8
$ws.Range("W5").Value = 2.87
$ws.Range("X5").Value = 1.31
$ws.Range("Y5").Value = 1.4
$ws.Range("Z5").Value = 2.55
$ws.Range("AA5").Value = 1.62
$ws.Range("AB5").Value = 2.02
$ws.Range("AG5").Value = 22
$ws.Range("AI5").Value = 10
$ws.Range("AJ5").Value = 6.2
$ws.Range("AL5").Value = 55
$ws.Range("AM5").Value = 8.75
$ws.Range("AN5").Value = 13
$ws.Range("AO5").Value = 9.5
$ws.Range("AQ5").Value = 20
$ws.Range("AR5").Value = 27
$ws.Range("AS5").Value = 400
$ws.Range("G6").Value = 1.73
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 2.3
$ws.Range("K6").Value = 2.25
$ws.Range("L6").Value = 4.33
$ws.Range("AC6").Value = 8
$ws.Range("AF6").Value = 15
$ws.Range("AI6").Value = 12
$ws.Range("AM6").Value = 13
$ws.Range("G7").Value = 2.5
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.7
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 3.4
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 7.5
$ws.Range("S7").Value = 2.1
$ws.Range("T7").Value = 1.7
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 11
$ws.Range("AI7").Value = 8.5
$ws.Range("AJ7").Value = 6
$ws.Range("AM7").Value = 8.5
$ws.Range("AN7").Value = 13
$ws.Range("AP7").Value = 26
$ws.Range("AS7").Value = 800
$ws.Range("G8").Value = 1.73
$ws.Range("H8").Value = 3.75
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 2.3
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 12
$ws.Range("AC8").Value = 8
$ws.Range("AD8").Value = 9
$ws.Range("AE8").Value = 8.5
$ws.Range("AF8").Value = 13
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 12
$ws.Range("AK8").Value = 15
$ws.Range("AM8").Value = 13
